$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 - Pabii blog post update
$ws.Range("D9").Value = "Global MBA 만든 뒷 이야기 – 3. 음해, 협잡의 근본적인 퇴출은 다수의 실력이 올라갔을 때"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/global-mba-behind-story-3/#utm_source=rss&utm_medium=rss&utm_campaign=global-mba-behind-story-3"

# Row 27 - ScatterLab tech blog post update
$ws.Range("D27").Value = "Apache Beam으로 머신러닝 데이터 파이프라인 구축하기 3편 - RunInference로 모델 추론하기"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/apache-beam-3/"

# Row 28 - Ropiens tistory post update
$ws.Range("D28").Value = "[논문리뷰]Sold!: Auction methods for multirobot coordination"
$ws.Range("E28").Value = "https://ropiens.tistory.com/203"

# Row 32 - Dodonam tistory post update
$ws.Range("D32").Value = "PMI(Pointwise Mutual Information); 점별 상호 정보량"
$ws.Range("E32").Value = "https://dodonam.tistory.com/395"
